# Applies the diff: updates the comment for the 2014-09-28 entry (row 10),
# corrects its Stop time / Delta Time, and appends a new log entry (row 12)
# for 2014-09-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOGT")

# --- Row 10: 2014-09-28 entry -------------------------------------------
# Stop time corrected from 16:00 to 15:10 (Delta Time formula recalculates
# from 95 to 45 automatically).
$ws.Range("C10").Value = 0.631944444444444

# Comment text changed from "Cree el esquema del documento de
# requerimientos." to "Trabajé en la creación del esquema del documento de
# requerimientos."
$ws.Range("H10").Value = "Trabajé en la creación del esquema del documento de requerimientos."

# --- Row 12 (new): 2014-09-29 entry --------------------------------------
$ws.Range("A12").Value = 41911
$ws.Range("B12").Value = 0.340277777777778
$ws.Range("C12").Value = 0.395833333333333
$ws.Range("D12").Value = 8
$ws.Range("E12").Formula = "=((HOUR(C12)-HOUR(B12))*60)+(MINUTE(C12)-MINUTE(B12))-D12"
$ws.Range("F12").Value = 5
$ws.Range("H12").Value = "Terminé la creación del esquema del documento de requerimientos."

# Row heights: rows that now hold a 2-line wrapped comment grow to 26.65pt.
$ws.Rows.Item(10).RowHeight = 26.65
$ws.Rows.Item(12).RowHeight = 26.65
